$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: A1 "Sample ID" -> "Specimen_Number", F1 "T/N" -> "SAMPLE_TYPE"
$ws.Range("A1").Value = "Specimen_Number"
$ws.Range("F1").Value = "SAMPLE_TYPE"

# Give A1 its own bold-white-on-black centered header style (distinct from the
# other headers, which keep the existing teal style already applied to F1).
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Color = 16777215
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").Interior.ColorIndex = 1
$ws.Range("A1").Interior.PatternColorIndex = 1

# Move the active selection to F1, matching the saved view state.
$ws.Range("F1").Select() | Out-Null
